# Generate Report for Handback
# Both localized files (cd40b744-... and a853dcf2-...) have now been handed
# back. Swap the two file rows on every sheet (a853dcf2 now listed first),
# and update their status / timestamps to reflect the handback.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("A2").Value = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.md"
$ovw.Range("B2").Value = "Handed back: in sync with en-US"
$ovw.Range("C2").Value = "Handed back: in sync with en-US"
$ovw.Range("A3").Value = "cd40b744-c5c6-440f-8fed-96225982b5fb.md"
$ovw.Range("B3").Value = "Handed back: in sync with en-US"
$ovw.Range("C3").Value = "Handed back: in sync with en-US"

foreach ($h in $ovw.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.md"
    }
    if ($addr -eq '$A$3') {
        $h.TextToDisplay = "cd40b744-c5c6-440f-8fed-96225982b5fb.md"
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("A2").Value = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.md"
$zh.Range("B2").Value = "Handed back: in sync with en-US"
$zh.Range("C2").Value = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.804e25099a650cff288ff18a2562a03e0504de04.zh-cn.xlf"
$zh.Range("D2").Value = "2016-02-24 10:07:24"
$zh.Range("E2").Value = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.md"
$zh.Range("F2").Value = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.804e25099a650cff288ff18a2562a03e0504de04.zh-cn.xlf"
$zh.Range("G2").Value = "2016-02-24 10:08:11"
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = "cd40b744-c5c6-440f-8fed-96225982b5fb.md"
$zh.Range("B3").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "cd40b744-c5c6-440f-8fed-96225982b5fb.d29c726b6406af72937382d8654471b2feaf05c5.zh-cn.xlf"
$zh.Range("D3").Value = "2016-02-24 10:07:24"
$zh.Range("E3").Value = "cd40b744-c5c6-440f-8fed-96225982b5fb.md"
$zh.Range("F3").Value = "cd40b744-c5c6-440f-8fed-96225982b5fb.d29c726b6406af72937382d8654471b2feaf05c5.zh-cn.xlf"
$zh.Range("G3").Value = "2016-02-24 10:08:11"
$zh.Range("H3").Value = "Include"

foreach ($h in $zh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.md" }
    if ($addr -eq '$C$2') { $h.TextToDisplay = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.804e25099a650cff288ff18a2562a03e0504de04.zh-cn.xlf" }
    if ($addr -eq '$E$2') { $h.TextToDisplay = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.md" }
    if ($addr -eq '$F$2') { $h.TextToDisplay = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.804e25099a650cff288ff18a2562a03e0504de04.zh-cn.xlf" }
    if ($addr -eq '$A$3') { $h.TextToDisplay = "cd40b744-c5c6-440f-8fed-96225982b5fb.md" }
    if ($addr -eq '$C$3') { $h.TextToDisplay = "cd40b744-c5c6-440f-8fed-96225982b5fb.d29c726b6406af72937382d8654471b2feaf05c5.zh-cn.xlf" }
    if ($addr -eq '$E$3') { $h.TextToDisplay = "cd40b744-c5c6-440f-8fed-96225982b5fb.md" }
    if ($addr -eq '$F$3') { $h.TextToDisplay = "cd40b744-c5c6-440f-8fed-96225982b5fb.d29c726b6406af72937382d8654471b2feaf05c5.zh-cn.xlf" }
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("A2").Value = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.md"
$de.Range("B2").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.804e25099a650cff288ff18a2562a03e0504de04.de-de.xlf"
$de.Range("D2").Value = "2016-02-24 10:07:36"
$de.Range("E2").Value = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.md"
$de.Range("F2").Value = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.804e25099a650cff288ff18a2562a03e0504de04.de-de.xlf"
$de.Range("G2").Value = "2016-02-24 10:08:33"
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = "cd40b744-c5c6-440f-8fed-96225982b5fb.md"
$de.Range("B3").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "cd40b744-c5c6-440f-8fed-96225982b5fb.d29c726b6406af72937382d8654471b2feaf05c5.de-de.xlf"
$de.Range("D3").Value = "2016-02-24 10:07:36"
$de.Range("E3").Value = "cd40b744-c5c6-440f-8fed-96225982b5fb.md"
$de.Range("F3").Value = "cd40b744-c5c6-440f-8fed-96225982b5fb.d29c726b6406af72937382d8654471b2feaf05c5.de-de.xlf"
$de.Range("G3").Value = "2016-02-24 10:08:33"
$de.Range("H3").Value = "Include"

foreach ($h in $de.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.md" }
    if ($addr -eq '$C$2') { $h.TextToDisplay = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.804e25099a650cff288ff18a2562a03e0504de04.de-de.xlf" }
    if ($addr -eq '$E$2') { $h.TextToDisplay = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.md" }
    if ($addr -eq '$F$2') { $h.TextToDisplay = "a853dcf2-d25f-4c0b-8e52-d9a0c8b2c8f4.804e25099a650cff288ff18a2562a03e0504de04.de-de.xlf" }
    if ($addr -eq '$A$3') { $h.TextToDisplay = "cd40b744-c5c6-440f-8fed-96225982b5fb.md" }
    if ($addr -eq '$C$3') { $h.TextToDisplay = "cd40b744-c5c6-440f-8fed-96225982b5fb.d29c726b6406af72937382d8654471b2feaf05c5.de-de.xlf" }
    if ($addr -eq '$E$3') { $h.TextToDisplay = "cd40b744-c5c6-440f-8fed-96225982b5fb.md" }
    if ($addr -eq '$F$3') { $h.TextToDisplay = "cd40b744-c5c6-440f-8fed-96225982b5fb.d29c726b6406af72937382d8654471b2feaf05c5.de-de.xlf" }
}
